$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this shifts existing rows 36-144 down to 37-145
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record's data
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value = "Maule"
$ws.Cells.Item(36, 4).Value = 45250
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100101
$ws.Cells.Item(36, 8).Value = "Berries"
$ws.Cells.Item(36, 9).Value = 100101001
$ws.Cells.Item(36, 10).Value = "Arándano (blue)"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 50
$ws.Cells.Item(36, 14).Value = 6000
$ws.Cells.Item(36, 15).Value = 6000
$ws.Cells.Item(36, 16).Value = 6000
$ws.Cells.Item(36, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(36, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(36, 19).Value = 3000
$ws.Cells.Item(36, 20).Value = 2
